$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preparacion de la Prueba (rows 2-3): fill in actual start/end times ---
$ws.Range("B2").Value = 0.006944444444444444
$ws.Range("D2").Value = 0.7361111111111112
$ws.Range("E2").Value = 0.7423611111111111

$ws.Range("B3").Value = 0.4166666666666667
$ws.Range("D3").Value = 0.7430555555555555
$ws.Range("E3").Value = 0.7465277777777778

# --- Incremento table (rows 6-8): new tasks for the Pila package ---
$ws.Range("A6").Value = "Crear interfaz Pila"
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 0.003472222222222222
$ws.Range("E6").Value = 0.37777777777777777
$ws.Range("F6").Value = 0.3826388888888889

$ws.Range("A7").Value = "Implementar Pila Estatica"
$ws.Range("B7").Value = 30
$ws.Range("C7").Value = 24
$ws.Range("D7").Value = 0.017361111111111112
$ws.Range("E7").Value = 0.3847222222222222
$ws.Range("F7").Value = 0.40625

$ws.Range("A8").Value = "Implementar Pila Dinamica"
$ws.Range("B8").Value = 40
$ws.Range("C8").Value = 34
$ws.Range("D8").Value = 0.024305555555555556
$ws.Range("E8").Value = 0.425
$ws.Range("F8").Value = 0.4527777777777778

# Match the original file's formatting quirk: H6:H10 all pick up H7's
# (full-bordered) style, then fill in the logical-error counts (all zero).
$ws.Range("H7").Copy()
$ws.Range("H6:H10").PasteSpecial(-4122)
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0

# Force B16 to recompute against the final totals (engine otherwise keeps a
# stale cached value from mid-edit dependency churn).
$ws.Range("B16").Formula = $ws.Range("B16").Formula

# --- New blank, selected cell B28 (picked up formatting from D12) ---
$ws.Range("D12").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("B28").Select()
